$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, shifting existing rows 136:186 down to 137:187
$ws.Rows("136:136").Insert()

# Populate the newly inserted row 136 with the new weekly record
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 45027
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112021
$ws.Range("G136").Value = "Ají"
$ws.Range("H136").Value = "Americana (o)"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 220
$ws.Range("K136").Value = 21000
$ws.Range("L136").Value = 22000
$ws.Range("M136").Value = 21455
$ws.Range("N136").Value = "`$/caja 25 kilos"
$ws.Range("O136").Value = "Provincia de Limarí"
$ws.Range("P136").Value = 858
$ws.Range("Q136").Value = 25
$ws.Range("R136").Value = "Hortaliza"
